$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.695.64"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "2.419.02"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.22"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.34"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.495"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").Value = "2.420.92"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -7.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.67"
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("D14").Value = "2.869.94"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "67.887.25"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000163"
$ws.Range("E16").Value = "  -6.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "22.65"
$ws.Range("E17").Value = "  -6.21%  "
$ws.Range("D18").Value = "2.424.31"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.20"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  -5.75%  "
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.78"
$ws.Range("E24").Value = "  -6.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.83"
$ws.Range("E25").Value = "  -4.78%  "
$ws.Range("D26").Value = "2.550.97"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.56"
$ws.Range("E27").Value = "  -7.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.87"
$ws.Range("E29").Value = "  -8.78%  "
$ws.Range("D30").Value = "0.0₃0791"
$ws.Range("E30").Value = "  -8.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.93"
$ws.Range("E31").Value = "  -9.06%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "408.58"
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  -5.92%  "
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.88"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  -5.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.51"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.295"
$ws.Range("E41").Value = "  -5.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.20"
$ws.Range("E42").Value = "  -7.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.43"
$ws.Range("E43").Value = "  -8.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.04"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.16"
$ws.Range("E45").Value = "  -5.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("E46").Value = "  -8.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  -5.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0706"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.547"
$ws.Range("E49").Value = "  -3.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.462"
$ws.Range("E50").Value = "  -9.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0896"
$ws.Range("E51").Value = "  -2.36%  "
